# Update odds values in the daily Betfair Back/Lay sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Australian A-League Men / Melbourne City vs Macarthur FC)
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 5.5
$ws.Range("I2").Value = 5.7
$ws.Range("J2").Value = 4.3
$ws.Range("L2").Value = 1.34
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 1.44
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.42
$ws.Range("Y2").Value = 21
$ws.Range("Z2").Value = 46
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 22
$ws.Range("AK2").Value = 16.5

# Row 3 (Algerian Ligue 1 / MC Alger vs ES Ben Aknoun)
$ws.Range("F3").Value = 1.38
$ws.Range("K3").Value = 5
$ws.Range("U3").Value = 1.5
$ws.Range("W3").Value = 3.25

# Row 4 (Portuguese Primeira Liga / Guimaraes vs Sporting Lisbon)
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 1.4
$ws.Range("I4").Value = 1.42
$ws.Range("K4").Value = 5.5
$ws.Range("N4").Value = 4.1
$ws.Range("Q4").Value = 1.88
$ws.Range("T4").Value = 2.14
$ws.Range("V4").Value = 3.35
$ws.Range("Z4").Value = 7.6
$ws.Range("AA4").Value = 11
$ws.Range("AF4").Value = 90
$ws.Range("AJ4").Value = 400
